$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("lolo", "dellomos", "09", "m", "manila", "6:15 A.M.", "26/07/2023", "None", "No, Fever, Cough"),
    @("kate", "dellomos", "s", "s", "s", "s", "26/07/2023", "First Dose", "Fever"),
    @("jade", "dellomos", "09", "s", "s", "s", "s", "Second Dose", "Sore Throat"),
    @("m", "dellomos", "m", "m", "m", "m", "m", "First Booster Shot", "Fever")
)

$startRow = 3
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $col = $j + 1
        $value = $rowData[$j]
        $cell = $ws.Cells.Item($row, $col)
        # Preserve values such as "09" as literal text instead of letting
        # Excel auto-convert them to numbers (which would drop the
        # leading zero).
        if ($value -match '^0[0-9]+$') {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $value
    }
}
